$wb = $excel.ActiveWorkbook

# ----- PIR sheet: append rows 164-176 -----
$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Range("A164:A176").NumberFormat = "@"
$wsPIR.Cells.Item(164,1).Value = "2026-01-28"
$wsPIR.Cells.Item(164,2).Value = "12:03:51"
$wsPIR.Cells.Item(164,3).Value = "12:00"
$wsPIR.Cells.Item(164,4).Value = "Bathroom"
$wsPIR.Cells.Item(164,5).Value = "No Motion"
$wsPIR.Cells.Item(164,6).Value = "Inactive"

$wsPIR.Cells.Item(165,1).Value = "2026-01-28"
$wsPIR.Cells.Item(165,2).Value = "12:03:55"
$wsPIR.Cells.Item(165,3).Value = "12:00"
$wsPIR.Cells.Item(165,4).Value = "Bathroom"
$wsPIR.Cells.Item(165,5).Value = "No Motion"
$wsPIR.Cells.Item(165,6).Value = "Inactive"

$wsPIR.Cells.Item(166,1).Value = "2026-01-28"
$wsPIR.Cells.Item(166,2).Value = "12:04:00"
$wsPIR.Cells.Item(166,3).Value = "12:00"
$wsPIR.Cells.Item(166,4).Value = "Bathroom"
$wsPIR.Cells.Item(166,5).Value = "No Motion"
$wsPIR.Cells.Item(166,6).Value = "Inactive"

$wsPIR.Cells.Item(167,1).Value = "2026-01-28"
$wsPIR.Cells.Item(167,2).Value = "12:04:05"
$wsPIR.Cells.Item(167,3).Value = "12:00"
$wsPIR.Cells.Item(167,4).Value = "Bathroom"
$wsPIR.Cells.Item(167,5).Value = "No Motion"
$wsPIR.Cells.Item(167,6).Value = "Inactive"

$wsPIR.Cells.Item(168,1).Value = "2026-01-28"
$wsPIR.Cells.Item(168,2).Value = "12:04:10"
$wsPIR.Cells.Item(168,3).Value = "12:00"
$wsPIR.Cells.Item(168,4).Value = "Bathroom"
$wsPIR.Cells.Item(168,5).Value = "No Motion"
$wsPIR.Cells.Item(168,6).Value = "Inactive"

$wsPIR.Cells.Item(169,1).Value = "2026-01-28"
$wsPIR.Cells.Item(169,2).Value = "12:04:16"
$wsPIR.Cells.Item(169,3).Value = "12:00"
$wsPIR.Cells.Item(169,4).Value = "Bathroom"
$wsPIR.Cells.Item(169,5).Value = "No Motion"
$wsPIR.Cells.Item(169,6).Value = "Inactive"

$wsPIR.Cells.Item(170,1).Value = "2026-01-28"
$wsPIR.Cells.Item(170,2).Value = "12:04:20"
$wsPIR.Cells.Item(170,3).Value = "12:00"
$wsPIR.Cells.Item(170,4).Value = "Bathroom"
$wsPIR.Cells.Item(170,5).Value = "No Motion"
$wsPIR.Cells.Item(170,6).Value = "Inactive"

$wsPIR.Cells.Item(171,1).Value = "2026-01-28"
$wsPIR.Cells.Item(171,2).Value = "12:04:25"
$wsPIR.Cells.Item(171,3).Value = "12:00"
$wsPIR.Cells.Item(171,4).Value = "Bathroom"
$wsPIR.Cells.Item(171,5).Value = "No Motion"
$wsPIR.Cells.Item(171,6).Value = "Inactive"

$wsPIR.Cells.Item(172,1).Value = "2026-01-28"
$wsPIR.Cells.Item(172,2).Value = "12:04:30"
$wsPIR.Cells.Item(172,3).Value = "12:00"
$wsPIR.Cells.Item(172,4).Value = "Bathroom"
$wsPIR.Cells.Item(172,5).Value = "No Motion"
$wsPIR.Cells.Item(172,6).Value = "Inactive"

$wsPIR.Cells.Item(173,1).Value = "2026-01-28"
$wsPIR.Cells.Item(173,2).Value = "12:04:36"
$wsPIR.Cells.Item(173,3).Value = "12:00"
$wsPIR.Cells.Item(173,4).Value = "Bathroom"
$wsPIR.Cells.Item(173,5).Value = "No Motion"
$wsPIR.Cells.Item(173,6).Value = "Inactive"

$wsPIR.Cells.Item(174,1).Value = "2026-01-28"
$wsPIR.Cells.Item(174,2).Value = "12:04:40"
$wsPIR.Cells.Item(174,3).Value = "12:00"
$wsPIR.Cells.Item(174,4).Value = "Bathroom"
$wsPIR.Cells.Item(174,5).Value = "No Motion"
$wsPIR.Cells.Item(174,6).Value = "Inactive"

$wsPIR.Cells.Item(175,1).Value = "2026-01-28"
$wsPIR.Cells.Item(175,2).Value = "12:04:45"
$wsPIR.Cells.Item(175,3).Value = "12:00"
$wsPIR.Cells.Item(175,4).Value = "Bathroom"
$wsPIR.Cells.Item(175,5).Value = "No Motion"
$wsPIR.Cells.Item(175,6).Value = "Inactive"

$wsPIR.Cells.Item(176,1).Value = "2026-01-28"
$wsPIR.Cells.Item(176,2).Value = "12:04:50"
$wsPIR.Cells.Item(176,3).Value = "12:00"
$wsPIR.Cells.Item(176,4).Value = "Bathroom"
$wsPIR.Cells.Item(176,5).Value = "No Motion"
$wsPIR.Cells.Item(176,6).Value = "Inactive"

# ----- Humidity sheet: append rows 154-165 -----
$wsHum = $wb.Worksheets.Item("Humidity")
$wsHum.Range("A154:A165").NumberFormat = "@"
$wsHum.Range("E154:E165").NumberFormat = "@"
$wsHum.Cells.Item(154,1).Value = "2026-01-28"
$wsHum.Cells.Item(154,2).Value = "12:03:51"
$wsHum.Cells.Item(154,3).Value = "12:00"
$wsHum.Cells.Item(154,4).Value = "Bathroom"
$wsHum.Cells.Item(154,5).Value = "88.0%"
$wsHum.Cells.Item(154,6).Value = "Active"

$wsHum.Cells.Item(155,1).Value = "2026-01-28"
$wsHum.Cells.Item(155,2).Value = "12:03:55"
$wsHum.Cells.Item(155,3).Value = "12:00"
$wsHum.Cells.Item(155,4).Value = "Bathroom"
$wsHum.Cells.Item(155,5).Value = "88.1%"
$wsHum.Cells.Item(155,6).Value = "Active"

$wsHum.Cells.Item(156,1).Value = "2026-01-28"
$wsHum.Cells.Item(156,2).Value = "12:03:59"
$wsHum.Cells.Item(156,3).Value = "12:00"
$wsHum.Cells.Item(156,4).Value = "Bathroom"
$wsHum.Cells.Item(156,5).Value = "87.2%"
$wsHum.Cells.Item(156,6).Value = "Active"

$wsHum.Cells.Item(157,1).Value = "2026-01-28"
$wsHum.Cells.Item(157,2).Value = "12:04:03"
$wsHum.Cells.Item(157,3).Value = "12:00"
$wsHum.Cells.Item(157,4).Value = "Bathroom"
$wsHum.Cells.Item(157,5).Value = "88.0%"
$wsHum.Cells.Item(157,6).Value = "Active"

$wsHum.Cells.Item(158,1).Value = "2026-01-28"
$wsHum.Cells.Item(158,2).Value = "12:04:07"
$wsHum.Cells.Item(158,3).Value = "12:00"
$wsHum.Cells.Item(158,4).Value = "Bathroom"
$wsHum.Cells.Item(158,5).Value = "87.0%"
$wsHum.Cells.Item(158,6).Value = "Active"

$wsHum.Cells.Item(159,1).Value = "2026-01-28"
$wsHum.Cells.Item(159,2).Value = "12:04:11"
$wsHum.Cells.Item(159,3).Value = "12:00"
$wsHum.Cells.Item(159,4).Value = "Bathroom"
$wsHum.Cells.Item(159,5).Value = "88.0%"
$wsHum.Cells.Item(159,6).Value = "Active"

$wsHum.Cells.Item(160,1).Value = "2026-01-28"
$wsHum.Cells.Item(160,2).Value = "12:04:15"
$wsHum.Cells.Item(160,3).Value = "12:00"
$wsHum.Cells.Item(160,4).Value = "Bathroom"
$wsHum.Cells.Item(160,5).Value = "88.0%"
$wsHum.Cells.Item(160,6).Value = "Active"

$wsHum.Cells.Item(161,1).Value = "2026-01-28"
$wsHum.Cells.Item(161,2).Value = "12:04:19"
$wsHum.Cells.Item(161,3).Value = "12:00"
$wsHum.Cells.Item(161,4).Value = "Bathroom"
$wsHum.Cells.Item(161,5).Value = "87.1%"
$wsHum.Cells.Item(161,6).Value = "Active"

$wsHum.Cells.Item(162,1).Value = "2026-01-28"
$wsHum.Cells.Item(162,2).Value = "12:04:23"
$wsHum.Cells.Item(162,3).Value = "12:00"
$wsHum.Cells.Item(162,4).Value = "Bathroom"
$wsHum.Cells.Item(162,5).Value = "88.1%"
$wsHum.Cells.Item(162,6).Value = "Active"

$wsHum.Cells.Item(163,1).Value = "2026-01-28"
$wsHum.Cells.Item(163,2).Value = "12:04:31"
$wsHum.Cells.Item(163,3).Value = "12:00"
$wsHum.Cells.Item(163,4).Value = "Bathroom"
$wsHum.Cells.Item(163,5).Value = "87.2%"
$wsHum.Cells.Item(163,6).Value = "Active"

$wsHum.Cells.Item(164,1).Value = "2026-01-28"
$wsHum.Cells.Item(164,2).Value = "12:04:35"
$wsHum.Cells.Item(164,3).Value = "12:00"
$wsHum.Cells.Item(164,4).Value = "Bathroom"
$wsHum.Cells.Item(164,5).Value = "88.1%"
$wsHum.Cells.Item(164,6).Value = "Active"

$wsHum.Cells.Item(165,1).Value = "2026-01-28"
$wsHum.Cells.Item(165,2).Value = "12:04:39"
$wsHum.Cells.Item(165,3).Value = "12:00"
$wsHum.Cells.Item(165,4).Value = "Bathroom"
$wsHum.Cells.Item(165,5).Value = "87.1%"
$wsHum.Cells.Item(165,6).Value = "Active"

# ----- Temperature sheet: append rows 154-165 -----
$wsTemp = $wb.Worksheets.Item("Temperature")
$wsTemp.Range("A154:A165").NumberFormat = "@"
$wsTemp.Cells.Item(154,1).Value = "2026-01-28"
$wsTemp.Cells.Item(154,2).Value = "12:03:52"
$wsTemp.Cells.Item(154,3).Value = "12:00"
$wsTemp.Cells.Item(154,4).Value = "Bathroom"
$wsTemp.Cells.Item(154,5).Value = "23.0C"
$wsTemp.Cells.Item(154,6).Value = "Active"

$wsTemp.Cells.Item(155,1).Value = "2026-01-28"
$wsTemp.Cells.Item(155,2).Value = "12:03:55"
$wsTemp.Cells.Item(155,3).Value = "12:00"
$wsTemp.Cells.Item(155,4).Value = "Bathroom"
$wsTemp.Cells.Item(155,5).Value = "23.0C"
$wsTemp.Cells.Item(155,6).Value = "Active"

$wsTemp.Cells.Item(156,1).Value = "2026-01-28"
$wsTemp.Cells.Item(156,2).Value = "12:03:59"
$wsTemp.Cells.Item(156,3).Value = "12:00"
$wsTemp.Cells.Item(156,4).Value = "Bathroom"
$wsTemp.Cells.Item(156,5).Value = "23.0C"
$wsTemp.Cells.Item(156,6).Value = "Active"

$wsTemp.Cells.Item(157,1).Value = "2026-01-28"
$wsTemp.Cells.Item(157,2).Value = "12:04:03"
$wsTemp.Cells.Item(157,3).Value = "12:00"
$wsTemp.Cells.Item(157,4).Value = "Bathroom"
$wsTemp.Cells.Item(157,5).Value = "23.0C"
$wsTemp.Cells.Item(157,6).Value = "Active"

$wsTemp.Cells.Item(158,1).Value = "2026-01-28"
$wsTemp.Cells.Item(158,2).Value = "12:04:07"
$wsTemp.Cells.Item(158,3).Value = "12:00"
$wsTemp.Cells.Item(158,4).Value = "Bathroom"
$wsTemp.Cells.Item(158,5).Value = "23.0C"
$wsTemp.Cells.Item(158,6).Value = "Active"

$wsTemp.Cells.Item(159,1).Value = "2026-01-28"
$wsTemp.Cells.Item(159,2).Value = "12:04:11"
$wsTemp.Cells.Item(159,3).Value = "12:00"
$wsTemp.Cells.Item(159,4).Value = "Bathroom"
$wsTemp.Cells.Item(159,5).Value = "23.0C"
$wsTemp.Cells.Item(159,6).Value = "Active"

$wsTemp.Cells.Item(160,1).Value = "2026-01-28"
$wsTemp.Cells.Item(160,2).Value = "12:04:15"
$wsTemp.Cells.Item(160,3).Value = "12:00"
$wsTemp.Cells.Item(160,4).Value = "Bathroom"
$wsTemp.Cells.Item(160,5).Value = "23.0C"
$wsTemp.Cells.Item(160,6).Value = "Active"

$wsTemp.Cells.Item(161,1).Value = "2026-01-28"
$wsTemp.Cells.Item(161,2).Value = "12:04:19"
$wsTemp.Cells.Item(161,3).Value = "12:00"
$wsTemp.Cells.Item(161,4).Value = "Bathroom"
$wsTemp.Cells.Item(161,5).Value = "23.0C"
$wsTemp.Cells.Item(161,6).Value = "Active"

$wsTemp.Cells.Item(162,1).Value = "2026-01-28"
$wsTemp.Cells.Item(162,2).Value = "12:04:23"
$wsTemp.Cells.Item(162,3).Value = "12:00"
$wsTemp.Cells.Item(162,4).Value = "Bathroom"
$wsTemp.Cells.Item(162,5).Value = "23.0C"
$wsTemp.Cells.Item(162,6).Value = "Active"

$wsTemp.Cells.Item(163,1).Value = "2026-01-28"
$wsTemp.Cells.Item(163,2).Value = "12:04:31"
$wsTemp.Cells.Item(163,3).Value = "12:00"
$wsTemp.Cells.Item(163,4).Value = "Bathroom"
$wsTemp.Cells.Item(163,5).Value = "23.0C"
$wsTemp.Cells.Item(163,6).Value = "Active"

$wsTemp.Cells.Item(164,1).Value = "2026-01-28"
$wsTemp.Cells.Item(164,2).Value = "12:04:35"
$wsTemp.Cells.Item(164,3).Value = "12:00"
$wsTemp.Cells.Item(164,4).Value = "Bathroom"
$wsTemp.Cells.Item(164,5).Value = "23.0C"
$wsTemp.Cells.Item(164,6).Value = "Active"

$wsTemp.Cells.Item(165,1).Value = "2026-01-28"
$wsTemp.Cells.Item(165,2).Value = "12:04:39"
$wsTemp.Cells.Item(165,3).Value = "12:00"
$wsTemp.Cells.Item(165,4).Value = "Bathroom"
$wsTemp.Cells.Item(165,5).Value = "22.9C"
$wsTemp.Cells.Item(165,6).Value = "Active"

